$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held a block of sample/placeholder "同行員工"(accompanying
# employee) rows (rows 2-8) plus one real row (row 9: 員工02 / 員工01 / 10).
# Feature-complete cleanup: drop all the placeholder rows so the sheet only
# keeps the header and a single, real example row. Deleting rows 2:8 shifts
# the former row 9 up to row 2, carrying its (unstyled) formatting with it.
$ws.Rows("2:8").Delete()

# Update that surviving example row to the final pairing: 員工01 <-> 員工09.
$ws.Range("A2").Value = "員工01"
$ws.Range("B2").Value = "員工09"

# C2 already holds the text value "10" carried over from the old row 9.

# The former "spare" columns (C onward) are no longer meant to be expanded.
$ws.Range("C1:XFD1").EntireColumn.Collapsed = $true
